{"js": "// Apply the \"Added many more features\" edits to the Clover Lady review.\n// Each entry: [exact existing text to find, replacement text].\nconst replacements = [\n  [\n    \"Play Clover Lady Free Today! Review & Ratings\",\n    \"Play Clover Lady Slot Free\"\n  ],\n  [\n    \"Impressive graphics and audio\",\n    \"Immersive graphics and design\"\n  ],\n  [\n    \"Wide selection of bonus features\",\n    \"Variety of betting options and volatilities\"\n  ],\n  [\n    \"Flexible betting range and automatic spins\",\n    \"Engaging bonus features with potential for big wins\"\n  ],\n  [\n    \"Selectable volatility levels\",\n    \"Vibrant symbols and visual animations\"\n  ],\n  [\n    \"Bonus feature can take a long time to trigger\",\n    \"No progressive jackpot feature\"\n  ],\n  [\n    \"Explore the enchanting forest with Clover Lady. Read the review, play for free, and discover bonus features, graphics, and design. Compatible on all devices.\",\n    \"Read our review of Clover Lady slot and play for free. Discover immersive graphics and engaging bonus features.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [findText, newText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Clover Lady review.\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2\n$wdReplaceAll = 2\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\nReplace-AllText \"Play Clover Lady Free Today! Review & Ratings\" \"Play Clover Lady Slot Free\"\nReplace-AllText \"Impressive graphics and audio\" \"Immersive graphics and design\"\nReplace-AllText \"Wide selection of bonus features\" \"Variety of betting options and volatilities\"\nReplace-AllText \"Flexible betting range and automatic spins\" \"Engaging bonus features with potential for big wins\"\nReplace-AllText \"Selectable volatility levels\" \"Vibrant symbols and visual animations\"\nReplace-AllText \"Bonus feature can take a long time to trigger\" \"No progressive jackpot feature\"\nReplace-AllText \"Explore the enchanting forest with Clover Lady. Read the review, play for free, and discover bonus features, graphics, and design. Compatible on all devices.\" \"Read our review of Clover Lady slot and play for free. Discover immersive graphics and engaging bonus features.\"\n"}
